# Auto-generated Excel COM-interop script applying the scraped diff to before.xlsx
# Updates per-row market/profit metrics (H, I, J, K, L, M, N columns) across the
# ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets, incl. a few cell
# additions/removals where a row gained or lost a trailing column value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3543.9744
$ws.Range("J138").Value = 4804.65
$ws.Range("L138").Value = 14413.95
$ws.Range("N138").Value = -24693.95
$ws.Range("H141").Value = 1394.5
$ws.Range("I141").Value = 1243.125
$ws.Range("K141").Value = 3729.375
$ws.Range("M141").Value = 1450.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2161.8125
$ws.Range("I2").Value = 1967.9286
$ws.Range("K2").Value = 1967.9286
$ws.Range("M2").Value = -1854.9286
$ws.Range("H45").Value = 4177.88
$ws.Range("I45").Value = 4150.6665
$ws.Range("K45").Value = 4150.6665
$ws.Range("M45").Value = -3773.6665
$ws.Range("H61").Value = 3033966
$ws.Range("J61").Value = 2299
$ws.Range("L61").Value = 2299
$ws.Range("N61").Value = -2723
$ws.Range("H98").Value = 80342.664
$ws.Range("J98").Value = 80342.664
$ws.Range("L98").Value = 80342.664
$ws.Range("N98").Value = -86332.664
$ws.Range("H116").Value = 2161.8125
$ws.Range("I116").Value = 1967.9286
$ws.Range("K116").Value = 1967.9286
$ws.Range("M116").Value = 326.0714
$ws.Range("H122").Value = 3120.2942
$ws.Range("I122").Value = 1568.625
$ws.Range("K122").Value = 4705.875
$ws.Range("M122").Value = -2255.875
$ws.Range("H136").Value = 3033966
$ws.Range("J136").Value = 2299
$ws.Range("L136").Value = 6897
$ws.Range("N136").Value = -11997
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2161.8125
$ws.Range("I3").Value = 1967.9286
$ws.Range("K3").Value = 1967.9286
$ws.Range("M3").Value = -1853.9286
$ws.Range("H20").Value = 2429.6
$ws.Range("I20").Value = 2722.8572
$ws.Range("J20").Value = 2173
$ws.Range("K20").Value = 2722.8572
$ws.Range("L20").Value = 2173
$ws.Range("M20").Value = -2475.8572
$ws.Range("N20").Value = -2667
$ws.Range("H86").Value = 2610.818
$ws.Range("I86").Value = 2737.5
$ws.Range("J86").Value = 2273
$ws.Range("K86").Value = 2737.5
$ws.Range("L86").Value = 2273
$ws.Range("M86").Value = -1614.5
$ws.Range("N86").Value = -4519
$ws.Range("H89").Value = 2610.818
$ws.Range("I89").Value = 2737.5
$ws.Range("J89").Value = 2273
$ws.Range("K89").Value = 13687.5
$ws.Range("L89").Value = 11365
$ws.Range("M89").Value = -8071.5
$ws.Range("N89").Value = -22597
$ws.Range("H105").Value = 1346.5172
$ws.Range("I105").Value = 1372.2593
$ws.Range("K105").Value = 1372.2593
$ws.Range("M105").Value = 374.7407000000001
$ws.Range("H107").Value = 5365.375
$ws.Range("I107").Value = 4703.4287
$ws.Range("K107").Value = 4703.4287
$ws.Range("M107").Value = -2783.4287
$ws.Range("H134").Value = 445970.94
$ws.Range("I134").Value = 518280.78
$ws.Range("K134").Value = 1554842.34
$ws.Range("M134").Value = -1552307.34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 753.9259
$ws.Range("I22").Value = 775.2353000000001
$ws.Range("K22").Value = 775.2353000000001
$ws.Range("M22").Value = -425.2353000000001
$ws.Range("H31").Value = 4878.1924
$ws.Range("J31").Value = 6241.5
$ws.Range("L31").Value = 6241.5
$ws.Range("N31").Value = -6831.5
$ws.Range("H34").Value = 4878.1924
$ws.Range("J34").Value = 6241.5
$ws.Range("L34").Value = 6241.5
$ws.Range("N34").Value = -6645.5
$ws.Range("H58").Value = 516609.38
$ws.Range("I58").Value = 618454.9
$ws.Range("K58").Value = 618454.9
$ws.Range("M58").Value = -618251.9
$ws.Range("H103").Value = 15920.5
$ws.Range("I103").Value = 15920.5
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 15920.5
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -14748.5
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 5564567.5
$ws.Range("I132").Value = 9384.535
$ws.Range("K132").Value = 28153.605
$ws.Range("M132").Value = -25623.605
$ws.Range("H134").Value = 2172.5625
$ws.Range("I134").Value = 2183.1333
$ws.Range("K134").Value = 6549.3999
$ws.Range("M134").Value = -4014.3999
$ws.Range("H136").Value = 516609.38
$ws.Range("I136").Value = 618454.9
$ws.Range("K136").Value = 1855364.7
$ws.Range("M136").Value = -1852814.7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 993.5
$ws.Range("I69").Value = 993.5
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 2980.5
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -2169.5
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 993.5
$ws.Range("I72").Value = 993.5
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 8941.5
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -4885.5
$ws.Range("N72").ClearContents()
$ws.Range("H81").Value = 6722.143
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 6722.143
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H93").Value = 3800.1
$ws.Range("J93").Value = 4714.5713
$ws.Range("L93").Value = 14143.7139
$ws.Range("N93").Value = -17887.7139
$ws.Range("H122").Value = 604741.9399999999
$ws.Range("J122").Value = 772574.75
$ws.Range("L122").Value = 6953172.75
$ws.Range("N122").Value = -6958072.75
$ws.Range("H131").Value = 12711.8
$ws.Range("I131").Value = 899
$ws.Range("J131").Value = 15665
$ws.Range("K131").Value = 2697
$ws.Range("L131").Value = 46995
$ws.Range("M131").Value = 2343
$ws.Range("N131").Value = -57075
$ws.Range("H132").Value = 3327.5
$ws.Range("I132").Value = 2198.5
$ws.Range("J132").Value = 4174.25
$ws.Range("K132").Value = 19786.5
$ws.Range("L132").Value = 37568.25
$ws.Range("M132").Value = -17256.5
$ws.Range("N132").Value = -42628.25
$ws.Range("H133").Value = 4136.1724
$ws.Range("J133").Value = 4999.0835
$ws.Range("L133").Value = 14997.2505
$ws.Range("N133").Value = -25117.2505
$ws.Range("H134").Value = 1154.4546
$ws.Range("I134").Value = 669.9
$ws.Range("K134").Value = 2009.7
$ws.Range("M134").Value = 3060.3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1864.1177
$ws.Range("I102").Value = 1358.7778
$ws.Range("K102").Value = 1358.7778
$ws.Range("M102").Value = 263.2221999999999
$ws.Range("H126").Value = 696832.6
$ws.Range("I126").Value = 2382398.2
$ws.Range("K126").Value = 7147194.600000001
$ws.Range("M126").Value = -7144724.600000001
$ws.Range("H132").Value = 755897
$ws.Range("I132").Value = 1097415.5
$ws.Range("K132").Value = 3292246.5
$ws.Range("M132").Value = -3289716.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3002.4
$ws.Range("I7").Value = 2753.3
$ws.Range("J7").Value = 3500.6
$ws.Range("K7").Value = 2753.3
$ws.Range("L7").Value = 3500.6
$ws.Range("M7").Value = -2641.3
$ws.Range("N7").Value = -3724.6
$ws.Range("H16").Value = 709.13336
$ws.Range("J16").Value = 875.75
$ws.Range("L16").Value = 875.75
$ws.Range("N16").Value = -1215.75
$ws.Range("H36").Value = 90000
$ws.Range("J36").Value = 90000
$ws.Range("L36").Value = 90000
$ws.Range("N36").Value = -91124
$ws.Range("H40").Value = 4205.222
$ws.Range("I40").Value = 3592.8635
$ws.Range("J40").Value = 6899.6
$ws.Range("K40").Value = 3592.8635
$ws.Range("L40").Value = 6899.6
$ws.Range("M40").Value = -3456.8635
$ws.Range("N40").Value = -7171.6
$ws.Range("H103").Value = 93598.39999999999
$ws.Range("J103").Value = 93598.39999999999
$ws.Range("L103").Value = 93598.39999999999
$ws.Range("N103").Value = -95942.39999999999
$ws.Range("H122").Value = 3542.5527
$ws.Range("I122").Value = 3268.037
$ws.Range("J122").Value = 4216.364
$ws.Range("K122").Value = 9804.110999999999
$ws.Range("L122").Value = 12649.092
$ws.Range("M122").Value = -7354.110999999999
$ws.Range("N122").Value = -17549.092
$ws.Range("H126").Value = 3002.4
$ws.Range("I126").Value = 2753.3
$ws.Range("J126").Value = 3500.6
$ws.Range("K126").Value = 8259.900000000001
$ws.Range("L126").Value = 10501.8
$ws.Range("M126").Value = -5789.900000000001
$ws.Range("N126").Value = -15441.8
$ws.Range("H132").Value = 1731402.8
$ws.Range("I132").Value = 1923364.1
$ws.Range("J132").Value = 3750
$ws.Range("K132").Value = 5770092.300000001
$ws.Range("L132").Value = 11250
$ws.Range("M132").Value = -5767562.300000001
$ws.Range("N132").Value = -16310
$ws.Range("H137").Value = 102135.4
$ws.Range("J137").Value = 115071.75
$ws.Range("L137").Value = 115071.75
$ws.Range("N137").Value = -125271.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4793969
$ws.Range("I132").Value = 4910529.5
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 14731588.5
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -14729058.5
$ws.Range("N132").Value = -50060
$ws.Range("H136").Value = 6671241
$ws.Range("I136").Value = 8448900
$ws.Range("K136").Value = 25346700
$ws.Range("M136").Value = -25344150
